$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StuffToImport")

# --- 1. Move the TopTable data block from I1:M4 down to A13:E16 --------
# Range.Cut preserves the original cell types (text vs number) and clears
# the source range in one go - exactly what a drag/cut-paste move does.
$ws.Range("I1:M4").Cut($ws.Range("A13"))

# Re-anchor the TopTable ListObject onto its new home range.
$loTop = $ws.ListObjects.Item("TopTable")
$loTop.Resize($ws.Range("A13:E16"))

# --- 2. Widen BaseTable with two new columns: Icon, Tooltip ------------
$loBase = $ws.ListObjects.Item("BaseTable")
$loBase.Resize($ws.Range("A1:G4"))
$ws.Range("F1").Value = "Icon"
$ws.Range("G1").Value = "Tooltip"

# --- 3. Widen TopTable (now at A13:E16) with the same two columns ------
$loTop.Resize($ws.Range("A13:G16"))
$ws.Range("F13").Value = "Icon"
$ws.Range("G13").Value = "Tooltip"

# --- 4. Sample data for the new Icon column (BaseTable row 2) ----------
$ws.Range("F2").Value = "test"

# --- 5. Tidy up the view: selection + phonetic settings ----------------
$ws.Range("G8").Select()
$ws.Range("A1:G16").SetPhonetic()
